$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row above the old row 15 ("phi / servoAngle" row), which
#    pushes it down to row 16 and shifts all formula references to C15 ->
#    C16 automatically (matches the diff: G5/G7 now use RADIANS(C16)).
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new "gripper" row (15). Borrow the border formatting from
#    the B6:D6 row (same border shape as the new row) before writing values.
# ---------------------------------------------------------------------------
$ws.Range("B6:D6").Copy()
$ws.Range("B15:D15").PasteSpecial(-4122)
$ws.Range("A1").Select()

$ws.Range("B15").Value = "gripper"
$ws.Range("C15").Value = 17
$ws.Range("D15").Value = "cm"

# ---------------------------------------------------------------------------
# 3. Row 16 is the old "phi" row (was row 15); update its value only - label
#    and units stay the same, just the servo angle test value changes.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 0

# ---------------------------------------------------------------------------
# 4. Update the leg-geometry input values (DEBUG_PRINT test numbers that
#    drive the servo angles to ~90 degrees).
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = -61.285
$ws.Range("C7").Value = 0
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 7.5
$ws.Range("C13").Value = 24
$ws.Range("C14").Value = 51.1

# ---------------------------------------------------------------------------
# 5. The G6 ("y") formula now also accounts for the gripper offset in C15.
# ---------------------------------------------------------------------------
$ws.Range("G6").Formula = "=C6+C15"

# ---------------------------------------------------------------------------
# 6. Add the new "ServoAngle" result row (13), pushing the highlighted
#    "final answer" box down from row 12 to row 13. Copy the existing
#    F12:H12 formatting (bold label + borders + yellow highlight) down to
#    F13:H13 first, then restore F12:H12 to the plain "mid-table" look
#    used by F10:H11.
# ---------------------------------------------------------------------------
$ws.Range("F12:H12").Copy()
$ws.Range("F13:H13").PasteSpecial(-4122)
$ws.Range("A1").Select()

$ws.Range("F11:H11").Copy()
$ws.Range("F12:H12").PasteSpecial(-4122)
$ws.Range("A1").Select()

$ws.Range("F12").Value = "Angle"
$ws.Range("F13").Value = "ServoAngle"
$ws.Range("G13").Formula = "=G12+90"
$ws.Range("H13").Value = "degrees"

# ---------------------------------------------------------------------------
# Recalculate, zoom to 125% and leave the selection on C6 (matches the new
# sheetView in the diff).
# ---------------------------------------------------------------------------
$excel.Calculate()
$ws.Range("C6").Select()
$excel.ActiveWindow.Zoom = 125
